$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Point of Contact" (column B) assignments to reflect database
# deployment work being handed to Derek.
$ws.Range("B8").Value  = "Derek, Walter"
$ws.Range("B11").Value = "Derek"
$ws.Range("B12").Value = "Derek"
$ws.Range("B13").Value = "Derek"
$ws.Range("B15").Value = "Derek"

$ws.Range("B19").Value = "Derek"
$ws.Range("B20").Value = "Derek"
$ws.Range("B21").Value = "Derek"
$ws.Range("B22").Value = "Derek"

# Shrink the highlighted duration bar on rows 11-13 (Lookup Tables, Foreign
# Keys, Constraints) and extend the bars on rows 20-21, re-using the
# existing green timeline-highlight format (copy/paste-format keeps the
# style pointed at the shared theme fill instead of minting a new RGB fill).
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D11:E13").PasteSpecial(-4122) | Out-Null
$ws.Range("G20:K20").PasteSpecial(-4122) | Out-Null
$ws.Range("F21:J21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the active selection to B9.
$ws.Range("B9").Select()
